# "blockquotes should be single-spaced"
#
# The "Block Text" style (w:styleId="BlockText") is used for blockquotes.
# It is based on "Body Text", which is double-spaced (w:line="480"
# w:lineRule="auto") with no space after paragraphs. Give BlockText its
# own explicit single-spacing (w:line="240" w:lineRule="auto") plus
# 12pt (240 twips) of space after each paragraph, so blockquotes no
# longer inherit the double spacing used by normal body paragraphs.

$d = $word.ActiveDocument
$style = $d.Styles("BlockText")

$style.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceSingle
$style.ParagraphFormat.SpaceAfter = 12       # 12 pt = 240 twips
